# ---------------------------------------------------------------------------
# "Avance de Fórmulas and Auditórias"
# Adds a new worksheet "CálculoAuditoríaV6" (an "Albarán" / delivery-note
# audit calculation) after the last existing sheet, fixes a shared-formula
# range bug on TrimestreONETWOv5, and leaves everything else untouched.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Fix the shared-formula ref bug on "TrimestreONETWOv5" (sheet D10 formula
#    was tagged with ref="D9:D12" but the master cell is D10, not D9).
#    Re-entering the formula on D10:D12 makes Excel recompute the correct
#    shared-formula anchor/ref ("D10:D12").
# ---------------------------------------------------------------------------
$wsTrim = $wb.Worksheets.Item("TrimestreONETWOv5")
$wsTrim.Range("D10:D12").Formula = "=B10+C10"

# ---------------------------------------------------------------------------
# 2) Add the new worksheet as the last tab and give it its name.
# ---------------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CálculoAuditoríaV6"

# Column widths (approximate "best fit" widths from the source file).
$ws.Columns.Item(1).ColumnWidth = 55.14
$ws.Columns.Item(2).ColumnWidth = 9.43
$ws.Columns.Item(3).ColumnWidth = 11.57

# ---------------------------------------------------------------------------
# 3) Header / info block (rows 1-4) + the Total/Descuento/IVA box (H2:I4).
# ---------------------------------------------------------------------------
$ws.Range("A1:G1").Value = "ALBARÁN"
$ws.Range("A2:G2").Value = "NOMBRE DE LA EMPRESA: ABMB S.A"
$ws.Range("H2").Value = "TOTAL"
$ws.Range("I2").Formula = "=SUM(F7:F21)"

$ws.Range("A3:G3").Value = "ACTIVIDAD: MATERIAL INFORMÁTICO Y ELECTRODOMÉSTICOS"
$ws.Range("H3").Value = "Descuento"
$ws.Range("I3").Value = 0.09

$ws.Range("A4:G4").Value = "NOMBRE DEL CLIENTE: SERVIPLUS  S.L."
$ws.Range("H4").Value = "I.V.A"
$ws.Range("I4").Value = 0.21

$ws.Range("I3:I4").NumberFormat = "0%"

# Row 1 styling: bottom-thin border under the title, shaded like the other
# sheet title banners.
$ws.Range("A1:G4").Font.Bold = $false

# Bottom border under row1 (title) and a separator line under row 5.
$r1 = $ws.Range("A1:G1").Borders.Item(9)
$r1.LineStyle = 1
$r1.Weight = 2

$r5 = $ws.Range("A5:G5").Borders.Item(9)
$r5.LineStyle = 1
$r5.Weight = 2

# Box border around the Total/Descuento/IVA summary (H2:I4).
$top = $ws.Range("H2:I2").Borders.Item(8)
$top.LineStyle = 1
$top.Weight = -4138
$left = $ws.Range("H2:H4").Borders.Item(7)
$left.LineStyle = 1
$left.Weight = -4138

# ---------------------------------------------------------------------------
# 4) Table header row (row 6) + data rows (7-21).
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "CONCEPTO"
$ws.Range("B6").Value = "UNIDADES"
$ws.Range("C6").Value = "PRE/UNIDAD"
$ws.Range("D6").Value = "DTO"
$ws.Range("E6").Value = "IVA"
$ws.Range("F6").Value = "TOTAL"

$items = @(
    @("Ordenador Pentium",       1,   1021.72),
    @("Impresora de Inyección",  2,   414.7),
    @("Monitor",                 5,   180.3),
    @("Televisor",               8,   570.96),
    @("Vídeo",                   8,   420.71),
    @("Diskettes",               860, 0.54),
    @("CD_ROM",                  9,   280.07),
    @("Tarjeta Controladora",    23,  33.66),
    @("Tarjeta VGA",             11,  28.25),
    @("Teclado",                 34,  32.45),
    @("Filtros de pantalla",     56,  66.11),
    @("Ratón",                   67,  23.44),
    @("Cable de impresora",      9,   5.89),
    @("Diskettera",              123, 31.25),
    @("Tarjeta de sonido",       74,  114.19)
)

$row = 7
foreach ($item in $items) {
    $ws.Range("A$row").Value = $item[0]
    $ws.Range("B$row").Value = $item[1]
    $ws.Range("C$row").Value = $item[2]
    $row++
}

# Formulas: D7/F7 standalone, then shared-formula blocks E7:E21, D8:D21,
# F8:F21 (matches the order/si grouping produced by the original edit).
$ws.Range("D7").Formula = "=C7*I`$3"
$ws.Range("E7:E21").Formula = "=C7*I`$4"
$ws.Range("F7").Formula = "=(C7-D7+E7)*B7"
$ws.Range("D8:D21").Formula = "=C8*I`$3"
$ws.Range("F8:F21").Formula = "=(C8-D8+E8)*B8"

# ---------------------------------------------------------------------------
# 5) Final view state: select G25 on the new sheet (matches the source).
# ---------------------------------------------------------------------------
$ws.Range("G25").Select()
